$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new LeetCode pandas entry (1517. Find Users With Valid E-Mails) as row 3

# Question (A3)
$ws.Range("A3").Value = "1517. Find Users With Valid E-Mails"

# Difficulty (B3): reuse the fill formatting used for "Easy" in row 2
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B3").PasteSpecial(-4122) | Out-Null
$ws.Range("B3").Value = "Easy"

# Pattern (C3)
$ws.Range("C3").Value = "String Methods"

# Notes (D3)
$ws.Range("D3").Value = "Use RegEx pattern with str.match(). When you use a dataframe twice, you filter the same dataframe based on the inner condition. E.g. users[users['mail'].str.match()]"

# Link (E3): reuse the Hyperlink cell style from row 2, fill the text, then attach the real hyperlink
$ws.Range("E2").Copy() | Out-Null
$ws.Range("E3").PasteSpecial(-4122) | Out-Null
$ws.Range("E3").Value = "https://leetcode.com/problems/find-users-with-valid-e-mails/solutions/3853585/regex-explained-pandas-mysql-an-effortless-and-simple-approach-with-comments/?envType=study-plan-v2&envId=30-days-of-pandas&lang=pythondata "
$ws.Hyperlinks.Add($ws.Range("E3"), "https://leetcode.com/problems/find-users-with-valid-e-mails/solutions/3853585/regex-explained-pandas-mysql-an-effortless-and-simple-approach-with-comments/?envType=study-plan-v2&envId=30-days-of-pandas&lang=pythondata ") | Out-Null
$ws.Range("E2").Copy() | Out-Null
$ws.Range("E3").PasteSpecial(-4122) | Out-Null

# Grow the table to cover the newly added row
$table = $ws.ListObjects.Item(1)
$table.Resize($ws.Range("A1:E3"))

# Match the author's final view state (scrolled right a column, selection on E15)
$ws.Range("E15").Select()
$excel.ActiveWindow.ScrollColumn = 2
